# "repro data for app" -- normalize the species presence columns:
#  - lower-case the header labels (drop the stray "y" shared string)
#  - turn the per-row "y" text markers into numeric 1 presence flags
#  - give columns B (month) and C (day) explicit widths
#  - leave the selection on G18, matching the author's last saved state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lower-case header labels for the species columns (D:G)
$ws.Range("D1").Value = "saccharina_latissima"
$ws.Range("E1").Value = "nereocystis_luetkeana"
$ws.Range("F1").Value = "alaria_marginata"
$ws.Range("G1").Value = "costaria_costata"

# Replace the "y" presence text markers with numeric 1 in D:G, rows 2-14
$presenceCells = "D2", "E3", "E5", "E6", "D7", "D8", "D9", "D10", "D11", "E11", "D12", "E12", "E13", "D14", "E14"
foreach ($addr in $presenceCells) {
    $ws.Range($addr).Value = 1
}

# Explicit column widths for B (month) and C (day)
$ws.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Restore the saved selection
$ws.Range("G18").Select()
